# Refresh the crypto symbol/price table (Sheet1) with the latest scrape.
# Price cells (column D) hold numeric-looking text (e.g. "250.60", "0.05677")
# that must stay as literal text (matches the source feed's formatting,
# including trailing zeros) rather than be auto-converted to a Number by
# Excel's input parser, so they are written with a leading apostrophe
# (the standard "force text" entry trick) before the digits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''250.60'
$ws.Range("D3").Value = '''21.80'
$ws.Range("D4").Value = '''5.580'
$ws.Range("D5").Value = '''0.05677'
$ws.Range("D6").Value = '''6.441'
$ws.Range("D7").Value = '''0.8080'
$ws.Range("D8").Value = '''1.039'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1427'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '''0.07234'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.03129'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.02918'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09272'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001667'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '''3.225'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04732'
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '''0.0005812'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("D18").Value = '''0.006458'
$ws.Range("D19").Value = '''0.005063'
$ws.Range("D22").Value = '''3.983'
$ws.Range("D23").Value = '''3.374'
$ws.Range("D24").Value = '''2.113'
$ws.Range("D25").Value = '''0.3294'
$ws.Range("D27").Value = '''0.0003101'
$ws.Range("D40").Value = '''0.04130'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1045'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("D42").Value = '''0.002970'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003246'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("D45").Value = '''0.00005647'
$ws.Range("D47").Value = '''0.7855'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range("D48").Value = '''0.01671'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("D49").Value = '''0.00002101'
